# Small modifications to the polling graphs
# Adds two new poll entries at the top of the "2019" sheet's table,
# pushing all existing rows down by two rows (formulas/shared refs
# shift accordingly), then fills the two new rows with the new poll
# data (Midgam/Channel 12[4] and Panel Project HaMidgam/Channel 13[5]).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# Duplicate the current rows 2:3 (which already carry the correct
# formulas/formatting for this table) and insert the copies above
# themselves - this shifts all the existing data down by two rows
# while keeping the per-row formulas (e.g. the V-column SUM) intact.
$ws.Rows("2:3").Copy()
$ws.Rows("2:3").Insert()

# New row 2: Midgam/Channel 12[4], 2019-04-10 (serial 43548)
$ws.Range("A2").Value = 43548
$ws.Range("B2").Value = "Midgam/Channel 12[4]"
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 32
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 28
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 7

# New row 3: Panel Project HaMidgam/Channel 13[5], 2019-04-10 (serial 43548)
$ws.Range("A3").Value = 43548
$ws.Range("B3").Value = "Panel Project HaMidgam/Channel 13[5]"
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 31
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 28
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 6

# Match the author's final cursor position recorded in the workbook.
$ws.Activate()
$ws.Range("M13").Select()
